# Weekly fruit/vegetable price update: a new daily record is inserted as a
# new row 15 (pushing the existing rows 15-69 down to 16-70), growing the
# sheet's used range from A1:R69 to A1:R70.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a brand-new row above the current row 15, shifting rows 15:69
# down to 16:70 (this also copies the date-number-format from the row
# above, which is what column D of this table relies on).
$ws.Rows(15).Insert()

# Populate the newly inserted row 15 with the new observation.
$ws.Range("A15").Value = 1
$ws.Range("B15").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C15").Value = "Arica y Parinacota"
$ws.Range("D15").Value = 44620
$ws.Range("E15").Value = 15
$ws.Range("F15").Value = 100112040
$ws.Range("G15").Value = "Cilantro"
$ws.Range("H15").Value = "Sin especificar"
$ws.Range("I15").Value = "Primera"
$ws.Range("J15").Value = 300
$ws.Range("K15").Value = 900
$ws.Range("L15").Value = 1000
$ws.Range("M15").Value = 950
$ws.Range("N15").Value = "$/atado 1,5 a 2 kilos"
$ws.Range("O15").Value = "Región de Arica y Parinacota"
$ws.Range("P15").Value = 475
$ws.Range("Q15").Value = 2
$ws.Range("R15").Value = "Hortaliza"
